$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Recalculated results (new data values)
$ws.Range("E7").Value = 0.1
$ws.Range("D10").Value = 0.45

# Move selection to D10 (also resets the scrolled-down "topLeftCell" view
# state back to showing the sheet from its top-left corner)
$ws.Activate()
$ws.Range("D10").Select() | Out-Null
